# Insert two new rows (new weekly records) before the existing row 601,
# pushing all subsequent rows (601-629) down to (603-631).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(601).Resize(2).Insert()

# New row 601
$ws.Range("A601").Value = 7
$ws.Range("B601").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C601").Value = "Ñuble"
$ws.Range("D601").Value = 44753
$ws.Range("E601").Value = 16
$ws.Range("F601").Value = "Fruta"
$ws.Range("G601").Value = 100108
$ws.Range("H601").Value = "Tropicales y subtropicales"
$ws.Range("I601").Value = 100108006
$ws.Range("J601").Value = "Plátano"
$ws.Range("K601").Value = "Sin especificar"
$ws.Range("L601").Value = 'Pintón'
$ws.Range("M601").Value = 80
$ws.Range("N601").Value = 22000
$ws.Range("O601").Value = 22000
$ws.Range("P601").Value = 22000
$ws.Range("Q601").Value = '$/caja 20 kilos'
$ws.Range("R601").Value = "Ecuador"
$ws.Range("S601").Value = 1100
$ws.Range("T601").Value = 20

# New row 602
$ws.Range("A602").Value = 7
$ws.Range("B602").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C602").Value = "Ñuble"
$ws.Range("D602").Value = 44753
$ws.Range("E602").Value = 16
$ws.Range("F602").Value = "Fruta"
$ws.Range("G602").Value = 100108
$ws.Range("H602").Value = "Tropicales y subtropicales"
$ws.Range("I602").Value = 100108006
$ws.Range("J602").Value = "Plátano"
$ws.Range("K602").Value = "Sin especificar"
$ws.Range("L602").Value = 'Primera Pintón'
$ws.Range("M602").Value = 160
$ws.Range("N602").Value = 23000
$ws.Range("O602").Value = 24000
$ws.Range("P602").Value = 23500
$ws.Range("Q602").Value = '$/caja 20 kilos'
$ws.Range("R602").Value = "Ecuador"
$ws.Range("S602").Value = 1175
$ws.Range("T602").Value = 20
